# Duplicate the existing flight rows (A6:H10) into a new block (A11:H15),
# as if the flight list was exported again / re-appended below the
# original data (see commit message: "export flight to pdf and excel").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$src = $ws.Range("A6:H10")
$dst = $ws.Range("A11:H15")
$src.Copy($dst)
